# Weekly update: a new price record for "Berenjena" @ Vega Modelo de Temuco
# was added for the current week. This shifts the old rows 226-239 down to
# 227-240 and inserts the new observation at row 226.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 226, pushing existing rows 226:239 down to 227:240
# (mirrors the rows-shifted-down pattern seen in the target diff).
$ws.Rows("226:226").Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(226, 1).Value  = 10
$ws.Cells.Item(226, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(226, 3).Value  = 'La Araucanía'
$ws.Cells.Item(226, 4).Value  = 44610
$ws.Cells.Item(226, 5).Value  = 9
$ws.Cells.Item(226, 6).Value  = 100112001
$ws.Cells.Item(226, 7).Value  = 'Berenjena'
$ws.Cells.Item(226, 8).Value  = 'Sin especificar'
$ws.Cells.Item(226, 9).Value  = 'Primera'
$ws.Cells.Item(226, 10).Value = 40
$ws.Cells.Item(226, 11).Value = 12000
$ws.Cells.Item(226, 12).Value = 12000
$ws.Cells.Item(226, 13).Value = 12000
$ws.Cells.Item(226, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(226, 15).Value = 'Región del Maule'
$ws.Cells.Item(226, 16).Value = 200
$ws.Cells.Item(226, 17).Value = 60
$ws.Cells.Item(226, 18).Value = 'Hortaliza'
